# Generate Report for Handoff
# Update status text and timestamps, then autofit the Status columns.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status moved from "In Translation" to "Ready for handoff" everywhere it appears.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Refresh generation / handoff timestamps.
$wsOverview.Range("G2").Value = "2016-09-01 02:47:41"
$wsZhCn.Range("H2").Value = "2016-09-01 02:47:37"
$wsDeDe.Range("H2").Value = "2016-09-01 02:47:41"

# Column widths adjust (auto-fit) to fit the new, longer status text.
$statusColumnWidth = 16.3333333
$wsOverview.Range("E1").ColumnWidth = $statusColumnWidth
$wsOverview.Range("F1").ColumnWidth = $statusColumnWidth
$wsZhCn.Range("C1").ColumnWidth = $statusColumnWidth
$wsDeDe.Range("C1").ColumnWidth = $statusColumnWidth
